$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 820, shifting existing rows 820:861 down to 821:862
$ws.Rows.Item(820).Insert()

# Populate the newly inserted row 820 with the new data point
$ws.Cells.Item(820, 1).NumberFormat = "@"
$ws.Cells.Item(820, 1).Value = "2026/02/20"
$ws.Cells.Item(820, 1).ClearFormats()
$ws.Cells.Item(820, 2).Value = "金"
$ws.Cells.Item(820, 3).Value = 13
$ws.Cells.Item(820, 4).Value = 54
